# Refresh the cryptos price/volume(1h) data (GitHub Actions scheduled update),
# and restore the correct Maker/EnergySwap row order (rows 50-51 were swapped).
#
# Column D ("Price") holds values as plain text in the source data, even when
# they look numeric (e.g. "147.04"). Excel's Range.Value setter auto-detects
# numeric-looking strings and silently coerces them to a Double (losing the
# text type and risking float rounding, e.g. 147.04 -> 147.03999999999999).
# A leading apostrophe is the standard Excel "force text" prefix: it is
# stripped from the stored value but keeps the cell as Text, matching the
# original inline-string cells. It's only needed for the D-column values
# that actually parse as numbers; everything else (links, names, and the
# "  +x.xx%  " volume strings, which always contain non-numeric characters)
# is already kept as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.365.50"
$ws.Range("E2").Value = "  +6.14%  "
$ws.Range("D3").Value = "2.488.80"
$ws.Range("E3").Value = "  +4.62%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'489.97"
$ws.Range("E5").Value = "  +7.16%  "
$ws.Range("D6").Value = "'147.04"
$ws.Range("E6").Value = "  +14.08%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  +7.51%  "
$ws.Range("D9").Value = "2.508.67"
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("D10").Value = "'5.80"
$ws.Range("E10").Value = "  +10.79%  "
$ws.Range("D11").Value = "'0.0984"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("D12").Value = "'0.334"
$ws.Range("E12").Value = "  +7.71%  "
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "2.918.83"
$ws.Range("E14").Value = "  +4.46%  "
$ws.Range("D15").Value = "56.348.26"
$ws.Range("E15").Value = "  +5.92%  "
$ws.Range("D16").Value = "'21.27"
$ws.Range("E16").Value = "  +9.56%  "
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  +6.13%  "
$ws.Range("D18").Value = "2.501.32"
$ws.Range("E18").Value = "  +5.02%  "
$ws.Range("D19").Value = "'4.57"
$ws.Range("E19").Value = "  +10.87%  "
$ws.Range("D20").Value = "'10.14"
$ws.Range("E20").Value = "  +9.43%  "
$ws.Range("D21").Value = "'319.32"
$ws.Range("E21").Value = "  +4.83%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'5.84"
$ws.Range("E23").Value = "  +10.66%  "
$ws.Range("D24").Value = "'58.73"
$ws.Range("E24").Value = "  +5.99%  "
$ws.Range("D25").Value = "'0.412"
$ws.Range("E25").Value = "  +8.32%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  +7.84%  "
$ws.Range("D28").Value = "2.589.30"
$ws.Range("E28").Value = "  +5.02%  "
$ws.Range("D29").Value = "'7.65"
$ws.Range("E29").Value = "  +9.36%  "
$ws.Range("D30").Value = "0.0₃0791"
$ws.Range("E30").Value = "  +10.49%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "'149.24"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").Value = "'18.26"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("E34").Value = "  +7.20%  "
$ws.Range("D35").Value = "'5.24"
$ws.Range("E35").Value = "  +5.85%  "
$ws.Range("E36").Value = "  +9.54%  "
$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = "  +7.32%  "
$ws.Range("D38").Value = "'0.862"
$ws.Range("E38").Value = "  +9.01%  "
$ws.Range("D39").Value = "'34.22"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("E40").Value = "  +8.79%  "
$ws.Range("D41").Value = "'0.0562"
$ws.Range("E41").Value = "  +8.20%  "
$ws.Range("D42").Value = "'0.612"
$ws.Range("E42").Value = "  +4.12%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  +9.11%  "
$ws.Range("D45").Value = "'4.81"
$ws.Range("E45").Value = "  +16.12%  "
$ws.Range("D46").Value = "'0.0922"
$ws.Range("E46").Value = "  +7.50%  "
$ws.Range("D47").Value = "'259.10"
$ws.Range("E47").Value = "  +19.37%  "
$ws.Range("E48").Value = "  +6.31%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'17.68"
$ws.Range("E50").Value = "  +8.45%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.896.28"
$ws.Range("E51").Value = "  -1.42%  "
